# Adds a "2022-Q4" quarter to the "600352-浙江龙盛" fund-holdings workbook:
#   1. A brand new worksheet "2022-Q4" (fund-level detail) is inserted right
#      after the "总计" (summary) sheet and before "2022-Q3".
#   2. The "总计" summary sheet gets a new data row for 2022-Q4 inserted at
#      the top of its table (row 2), pushing the rest of the quarters down
#      by one row; a fresh row for 2020-Q4 appears at the bottom (row 10).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT (no numeric coercion),
# matching the workbook's existing convention of storing these particular
# figures (fund codes, percentages, NAV amounts) as text rather than numbers,
# while leaving the cell's style/format untouched afterwards.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Helper: copy the formatting (style) of $srcCell onto $dstRange.
# ---------------------------------------------------------------------------
function Copy-Format {
    param($srcCell, $dstRange)
    $srcCell.Copy()
    $dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ===========================================================================
# 1. Insert the new "2022-Q4" worksheet right after "总计"
# ===========================================================================
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# NB: fetch this reference only *after* Worksheets.Add() above - a reference
# obtained before the collection is mutated ends up stale and silently
# fails to carry formats over via Copy/PasteSpecial below.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Match page margins used by every other sheet in the workbook.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Borrow the header/index style (bold + thin border + centered) that's
# already used for row 1 and column A on every other quarterly sheet.
$styleSrc = $q3Sheet.Range("B1")
Copy-Format $styleSrc $newSheet.Range("B1:H1")
Copy-Format $q3Sheet.Range("A2") $newSheet.Range("A2:A4")

# Row 1 - headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "090007"
$newSheet.Range("C2").Value = "大成策略回报混合"
Set-TextValue $newSheet.Range("D2") "14.14"
Set-TextValue $newSheet.Range("E2") "60.14"
Set-TextValue $newSheet.Range("F2") "2.19"
Set-TextValue $newSheet.Range("G2") "0.3097"
$newSheet.Range("H2").Value = 10

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "003318"
$newSheet.Range("C3").Value = "景顺长城中证500行业中性低波动指数"
Set-TextValue $newSheet.Range("D3") "10.01"
Set-TextValue $newSheet.Range("E3") "93.81"
Set-TextValue $newSheet.Range("F3") "1.10"
Set-TextValue $newSheet.Range("G3") "0.1101"
$newSheet.Range("H3").Value = 3

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "512260"
$newSheet.Range("C4").Value = "华安中证500行业中性低波动ETF"
Set-TextValue $newSheet.Range("D4") "0.94"
Set-TextValue $newSheet.Range("E4") "97.66"
Set-TextValue $newSheet.Range("F4") "1.15"
Set-TextValue $newSheet.Range("G4") "0.0108"
$newSheet.Range("H4").Value = 3

# ===========================================================================
# 2. Update the "总计" summary sheet: insert the 2022-Q4 row at the top of
#    the table and push every other quarter down by one row.
# ===========================================================================
$rows = @(
    @{ Date = "2022-Q4"; Count = 3;  Value = 0.43 },
    @{ Date = "2022-Q3"; Count = 4;  Value = 0.89 },
    @{ Date = "2022-Q2"; Count = 6;  Value = 3.23 },
    @{ Date = "2022-Q1"; Count = 8;  Value = 4.11 },
    @{ Date = "2021-Q4"; Count = 9;  Value = 3.85 },
    @{ Date = "2021-Q3"; Count = 11; Value = 4.15 },
    @{ Date = "2021-Q2"; Count = 12; Value = 18.21 },
    @{ Date = "2021-Q1"; Count = 19; Value = 20.15 },
    @{ Date = "2020-Q4"; Count = 12; Value = 7.3 }
)

# Row 10 is brand new (the table previously ended at row 9), so column A's
# "bold index" style needs to be carried onto it explicitly.
Copy-Format $totalSheet.Range("A9") $totalSheet.Range("A10")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $rows[$i].Date
    $totalSheet.Cells.Item($r, 3).Value = $rows[$i].Count
    $totalSheet.Cells.Item($r, 4).Value = $rows[$i].Value
}
